$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Preserve the current (style-21) look of row 228 onto the brand-new
#    row 238 BEFORE row 228's own formatting gets changed below.
# ---------------------------------------------------------------------------
$ws.Range("A228:G228").Copy()
$ws.Range("A238:G238").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Re-format rows 228:231 (B,C,D,F,G) from the old look to the look used
#    by the rest of the table (same as e.g. row 87), and format the brand
#    new rows 232:237 the same way.
# ---------------------------------------------------------------------------
$ws.Range("A87:G87").Copy()
$ws.Range("A228:G231").PasteSpecial(-4122)
$ws.Range("A87:G87").Copy()
$ws.Range("A232:G237").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Fill in the data values for the new rows (232:238). Rows 228:231 keep
#    their existing values - only their format changed above.
# ---------------------------------------------------------------------------
$data = @(
    @(232, 44109, 432, 3, 2336, 32, 8),
    @(233, 44110, 691, 3, 2936, 31, 8),
    @(234, 44111, 489, 3, 3351, 40, 13),
    @(235, 44112, 375, 4, 3703, 60, 20),
    @(236, 44113, 354, 2, 3863, 68, 25),
    @(237, 44114, 374, 2, 4161, 73, 28),
    @(238, 44115, 561, 8, 4587, 90, 29)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 5).Formula = "=D$r-F$r"
}
